$d = $word.ActiveDocument

# Update the date heading
$d.Content.Find.Execute("2024-04-22 Monday", $true, $false, $false, $false, $false, $true, 1, $false, "2024-04-23 Tuesday", 2) | Out-Null

# Update the division-table cell values
$t = $d.Tables.Item(1)
# Row 1
$t.Cell(1, 1).Range.Text = "44÷9=4, 8"
$t.Cell(1, 2).Range.Text = "88÷3=29, 1"
$t.Cell(1, 3).Range.Text = "69÷9=7, 6"
$t.Cell(1, 4).Range.Text = "67÷2=33, 1"
$t.Cell(1, 5).Range.Text = "19÷2=9, 1"

# Row 5
$t.Cell(5, 1).Range.Text = "15÷6=2, 3"
$t.Cell(5, 2).Range.Text = "60÷9=6, 6"
$t.Cell(5, 3).Range.Text = "56÷4=14, 0"
$t.Cell(5, 4).Range.Text = "85÷4=21, 1"
$t.Cell(5, 5).Range.Text = "26÷9=2, 8"

# Row 9
$t.Cell(9, 1).Range.Text = "12÷2=6, 0"
$t.Cell(9, 2).Range.Text = "29÷3=9, 2"
$t.Cell(9, 3).Range.Text = "56÷4=14, 0"
$t.Cell(9, 4).Range.Text = "41÷4=10, 1"
$t.Cell(9, 5).Range.Text = "33÷5=6, 3"

# Row 13
$t.Cell(13, 1).Range.Text = "89÷6=14, 5"
$t.Cell(13, 2).Range.Text = "64÷9=7, 1"
$t.Cell(13, 3).Range.Text = "31÷2=15, 1"
$t.Cell(13, 4).Range.Text = "41÷4=10, 1"
$t.Cell(13, 5).Range.Text = "41÷7=5, 6"

# Row 17
$t.Cell(17, 1).Range.Text = "23÷9=2, 5"
$t.Cell(17, 2).Range.Text = "63÷7=9, 0"
$t.Cell(17, 3).Range.Text = "58÷9=6, 4"
$t.Cell(17, 4).Range.Text = "67÷7=9, 4"
$t.Cell(17, 5).Range.Text = "38÷3=12, 2"
